$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data

$ws.Range("D2").Value = '42.766.55'
$ws.Range("E2").Value = '  -5.52%  '

$ws.Range("D3").Value = '2.212.02'
$ws.Range("E3").Value = '  -6.66%  '

$ws.Range("E4").Value = '  +0.02%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'314.29"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +1.20%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'98.09"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -10.00%  '

$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.580"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = '  -8.01%  '

$ws.Range("E8").Value = '  +0.00%  '

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.558"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  -9.63%  '

$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'36.44"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = '  -11.57%  '

$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'54.16"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  -2.32%  '

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.0823"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -10.46%  '

$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'7.74"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  -8.89%  '

$ws.Range("E14").Value = '  -3.85%  '

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'0.861"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -12.25%  '

$ws.Range("D16").Value = '2.547.73'
$ws.Range("E16").Value = '  -6.76%  '

$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'14.09"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = '  -7.86%  '

$ws.Range("D18").Value = '2.202.05'
$ws.Range("E18").Value = '  -7.15%  '

$ws.Range("D19").Value = '42.622.55'
$ws.Range("E19").Value = '  -5.86%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'14.70"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -1.88%  '

$ws.Range("D21").Value = '0.0₃0957'
$ws.Range("E21").Value = '  -9.82%  '

$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'6.37"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -12.65%  '

$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'64.95"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -11.41%  '

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'3.15"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -9.91%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'236.28"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  -9.14%  '

$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'2.12"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  -8.54%  '

$ws.Range("E27").Value = '  -0.03%  '

$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'10.02"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -10.14%  '

$ws.Range("E29").Value = '  -5.58%  '

$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'6.27"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = '  -13.78%  '

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'20.38"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  -9.15%  '

$ws.Range("E32").Value = '  -9.77%  '

$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'33.90"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  -10.19%  '

$style = $ws.Range("D34").Style
$ws.Range("D34").Value = "'155.06"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  -8.46%  '

$ws.Range("E35").Value = '  -6.32%  '

$ws.Range("E36").Value = '  +6.63%  '

$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'1.98"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +12.49%  '

$ws.Range("E38").Value = '  -6.72%  '

$ws.Range("E39").Value = '  -6.80%  '

$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.101"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -13.13%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'3.69"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -6.13%  '

$ws.Range("E42").Value = '  -8.46%  '

$ws.Range("D43").Value = '1.871.12'
$ws.Range("E43").Value = '  +6.31%  '

$ws.Range("E44").Value = '  +0.13%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'89.18"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -10.54%  '

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'12.09"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -6.54%  '

$ws.Range("E47").Value = '  -10.22%  '

$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'5.40"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -2.91%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'60.18"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  -13.38%  '

$ws.Range("B50").Value = 'ordi'
$ws.Range("C50").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'75.56"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  -7.34%  '

$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'8.60"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  -6.22%  '
